$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout: A1 = 0.9296 (number), A2 = 0 (number)
# Target layout:
#   A1 = "Fleiss' Kappa" (string), B1 = "valor p" (string)
#   A2 = 0.9296 (number),          B2 = 0 (number)

# Move the existing numeric values down a row, into columns A and B of row 2.
$ws.Range("B2").Value = $ws.Range("A2").Value()
$ws.Range("A2").Value = $ws.Range("A1").Value()

# Write the new header labels into row 1.
$ws.Range("A1").Value = "Fleiss' Kappa"
$ws.Range("B1").Value = "valor p"
